$wb = $excel.ActiveWorkbook

# --- Trees sheet: add new row 5 (new tree "edfasd" for user "tang") ---
$trees = $wb.Worksheets.Item("Trees")
# Seed row 5 from row 2 (same user/date/flow/notes pattern) then overwrite the differing cells.
$trees.Range("A2:K2").Copy($trees.Range("A5:K5"))
$trees.Cells.Item(5, 2).Value = "edfasd"
$trees.Cells.Item(5, 3).Value = 2
$trees.Cells.Item(5, 4).Value = 4
$trees.Cells.Item(5, 8).Value = 39.3989
$trees.Cells.Item(5, 9).Value = -74.5145
# Start of Season Notes (J5) needs to be a literal empty-text value (quote-prefix trick),
# then clear the resulting formatting back to Normal so no stray number format sticks.
$trees.Cells.Item(5, 10).Value = "'"
$trees.Cells.Item(5, 10).Style = "Normal"

# --- Seasons sheet: add new row 6 (same new tree "edfasd", season 2021) ---
$seasons = $wb.Worksheets.Item("Seasons")
$seasons.Range("A2:H2").Copy($seasons.Range("A6:H6"))
$seasons.Cells.Item(6, 2).Value = "edfasd"
# Start of Season Notes (G6) -> literal empty text, same trick as above.
$seasons.Cells.Item(6, 7).Value = "'"
$seasons.Cells.Item(6, 7).Style = "Normal"
